$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $escaped = $val -replace '"', '""'
    $r.Formula = '="' + $escaped + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)
}

Set-TextValue "D2" '27.701.90'
Set-TextValue "E2" '  -0.75%  '
Set-TextValue "D3" '1.591.37'
Set-TextValue "E3" '  -2.39%  '
Set-TextValue "E4" '  +0.17%  '
Set-TextValue "D5" '208.34'
Set-TextValue "D6" '0.501'
Set-TextValue "E6" '  -2.87%  '
Set-TextValue "E7" '  +0.21%  '
Set-TextValue "D8" '22.31'
Set-TextValue "E8" '  -3.99%  '
Set-TextValue "E9" '  -1.91%  '
Set-TextValue "E10" '  -2.34%  '
Set-TextValue "D11" '0.0869'
Set-TextValue "E11" '  -1.59%  '
Set-TextValue "D12" '1.816.37'
Set-TextValue "E12" '  -2.43%  '
Set-TextValue "D13" '1.587.71'
Set-TextValue "E13" '  -2.37%  '
Set-TextValue "E14" '  -3.75%  '
Set-TextValue "D15" '0.532'
Set-TextValue "E15" '  -4.28%  '
Set-TextValue "D16" '27.675.54'
Set-TextValue "E16" '  -0.88%  '
Set-TextValue "D17" '63.40'
Set-TextValue "E17" '  -2.16%  '
Set-TextValue "D18" '220.14'
Set-TextValue "E18" '  -3.39%  '
Set-TextValue "D19" '0.0₃0697'
Set-TextValue "E19" '  -2.99%  '
Set-TextValue "D20" '7.35'
Set-TextValue "E20" '  -3.39%  '
Set-TextValue "E21" '  +0.16%  '
Set-TextValue "D22" '4.15'
Set-TextValue "E22" '  -4.55%  '
Set-TextValue "D23" '9.70'
Set-TextValue "E23" '  -2.93%  '
Set-TextValue "E24" '  -3.89%  '
Set-TextValue "D25" '153.86'
Set-TextValue "E25" '  -0.56%  '
Set-TextValue "D26" '6.80'
Set-TextValue "E26" '  -1.76%  '
Set-TextValue "E27" '  +0.17%  '
Set-TextValue "D28" '15.15'
Set-TextValue "E28" '  -1.57%  '
Set-TextValue "E29" '  -4.91%  '
Set-TextValue "E30" '  -1.53%  '
Set-TextValue "E31" '  -2.19%  '
Set-TextValue "D32" '3.23'
Set-TextValue "E32" '  -4.87%  '
Set-TextValue "D33" '1.375.21'
Set-TextValue "E33" '  -2.94%  '
Set-TextValue "D34" '2.95'
Set-TextValue "E34" '  -5.01%  '
Set-TextValue "E35" '  -4.70%  '
Set-TextValue "D36" '0.974'
Set-TextValue "E36" '  -2.80%  '
Set-TextValue "D38" '0.0168'
Set-TextValue "E38" '  -1.27%  '
Set-TextValue "D39" '0.538'
Set-TextValue "D40" '0.830'
Set-TextValue "E40" '  -2.38%  '
Set-TextValue "E41" '  +0.22%  '
Set-TextValue "D42" '0.968'
Set-TextValue "E42" '  -3.70%  '
Set-TextValue "D43" '64.42'
Set-TextValue "E43" '  -2.05%  '
Set-TextValue "E44" '  +2.21%  '
Set-TextValue "E45" '  -3.62%  '
Set-TextValue "E46" '  -5.16%  '
Set-TextValue "D47" '1.727.83'
Set-TextValue "E47" '  -2.39%  '
Set-TextValue "D48" '86.99'
Set-TextValue "E48" '  -1.79%  '
Set-TextValue "E49" '  -1.04%  '
Set-TextValue "E50" '  -3.98%  '
Set-TextValue "E51" '  -1.46%  '

$excel.CutCopyMode = 0
